$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 16: A16 = date (serial 45978), B16 = 90, reusing the
# existing date number format from the cell above (A2) so no new
# style/numFmt entries get created.
$ws.Range("A2").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = [DateTime]::FromOADate(45978)
$ws.Range("B16").Value = 90

# Move the selection/active cell to D24 (just a cursor position change)
$ws.Range("D24").Select()
